# Re-do TODE grade norms with week coding for grade.
# Each worksheet (K-Fall, K-Spring, 1-Fall, 1-Spring, 2-Fall, 2-Spring) is a
# raw-score -> standard-score lookup table: column A holds the raw score,
# column B holds the corresponding standard score. This script rewrites the
# standard-score (column B) values per the updated norms; raw scores in
# column A and all other sheet structure are left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("K-Fall")
$ws.Cells.Item(2, 2).Value = 59
$ws.Cells.Item(3, 2).Value = 61
$ws.Cells.Item(4, 2).Value = 64
$ws.Cells.Item(5, 2).Value = 66
$ws.Cells.Item(6, 2).Value = 69
$ws.Cells.Item(7, 2).Value = 71
$ws.Cells.Item(8, 2).Value = 74
$ws.Cells.Item(9, 2).Value = 76
$ws.Cells.Item(10, 2).Value = 79
$ws.Cells.Item(11, 2).Value = 81
$ws.Cells.Item(12, 2).Value = 84
$ws.Cells.Item(13, 2).Value = 86
$ws.Cells.Item(14, 2).Value = 89
$ws.Cells.Item(15, 2).Value = 91
$ws.Cells.Item(16, 2).Value = 94
$ws.Cells.Item(17, 2).Value = 96
$ws.Cells.Item(18, 2).Value = 99
$ws.Cells.Item(19, 2).Value = 101
$ws.Cells.Item(20, 2).Value = 104
$ws.Cells.Item(21, 2).Value = 106
$ws.Cells.Item(22, 2).Value = 109
$ws.Cells.Item(23, 2).Value = 111
$ws.Cells.Item(24, 2).Value = 114
$ws.Cells.Item(25, 2).Value = 116
$ws.Cells.Item(26, 2).Value = 119
$ws.Cells.Item(27, 2).Value = 121
$ws.Cells.Item(28, 2).Value = 124
$ws.Cells.Item(29, 2).Value = 126

$ws = $wb.Worksheets.Item("K-Spring")
$ws.Cells.Item(2, 2).Value = 51
$ws.Cells.Item(3, 2).Value = 54
$ws.Cells.Item(4, 2).Value = 56
$ws.Cells.Item(5, 2).Value = 59
$ws.Cells.Item(6, 2).Value = 61
$ws.Cells.Item(7, 2).Value = 64
$ws.Cells.Item(8, 2).Value = 66
$ws.Cells.Item(9, 2).Value = 69
$ws.Cells.Item(10, 2).Value = 71
$ws.Cells.Item(11, 2).Value = 74
$ws.Cells.Item(12, 2).Value = 76
$ws.Cells.Item(13, 2).Value = 79
$ws.Cells.Item(14, 2).Value = 81
$ws.Cells.Item(15, 2).Value = 84
$ws.Cells.Item(16, 2).Value = 86
$ws.Cells.Item(17, 2).Value = 89
$ws.Cells.Item(18, 2).Value = 91
$ws.Cells.Item(21, 2).Value = 99
$ws.Cells.Item(22, 2).Value = 101
$ws.Cells.Item(23, 2).Value = 104
$ws.Cells.Item(24, 2).Value = 106
$ws.Cells.Item(25, 2).Value = 109
$ws.Cells.Item(26, 2).Value = 111
$ws.Cells.Item(27, 2).Value = 114
$ws.Cells.Item(28, 2).Value = 116
$ws.Cells.Item(29, 2).Value = 119
$ws.Cells.Item(30, 2).Value = 121
$ws.Cells.Item(31, 2).Value = 124
$ws.Cells.Item(32, 2).Value = 126
$ws.Cells.Item(33, 2).Value = 129

$ws = $wb.Worksheets.Item("1-Fall")
$ws.Cells.Item(2, 2).Value = 44
$ws.Cells.Item(3, 2).Value = 46
$ws.Cells.Item(4, 2).Value = 49
$ws.Cells.Item(5, 2).Value = 51
$ws.Cells.Item(6, 2).Value = 54
$ws.Cells.Item(7, 2).Value = 56
$ws.Cells.Item(8, 2).Value = 59
$ws.Cells.Item(9, 2).Value = 61
$ws.Cells.Item(10, 2).Value = 64
$ws.Cells.Item(11, 2).Value = 66
$ws.Cells.Item(12, 2).Value = 69
$ws.Cells.Item(13, 2).Value = 71
$ws.Cells.Item(14, 2).Value = 74
$ws.Cells.Item(15, 2).Value = 76
$ws.Cells.Item(16, 2).Value = 79
$ws.Cells.Item(17, 2).Value = 81
$ws.Cells.Item(18, 2).Value = 84
$ws.Cells.Item(20, 2).Value = 89
$ws.Cells.Item(21, 2).Value = 91
$ws.Cells.Item(22, 2).Value = 94
$ws.Cells.Item(23, 2).Value = 96
$ws.Cells.Item(24, 2).Value = 99
$ws.Cells.Item(25, 2).Value = 101
$ws.Cells.Item(26, 2).Value = 104
$ws.Cells.Item(27, 2).Value = 106
$ws.Cells.Item(28, 2).Value = 109
$ws.Cells.Item(29, 2).Value = 111
$ws.Cells.Item(30, 2).Value = 114
$ws.Cells.Item(31, 2).Value = 116
$ws.Cells.Item(32, 2).Value = 119
$ws.Cells.Item(33, 2).Value = 121
$ws.Cells.Item(34, 2).Value = 124

$ws = $wb.Worksheets.Item("1-Spring")
$ws.Cells.Item(2, 2).Value = 59
$ws.Cells.Item(3, 2).Value = 61
$ws.Cells.Item(4, 2).Value = 64
$ws.Cells.Item(5, 2).Value = 66
$ws.Cells.Item(6, 2).Value = 69
$ws.Cells.Item(7, 2).Value = 71
$ws.Cells.Item(8, 2).Value = 74
$ws.Cells.Item(9, 2).Value = 76
$ws.Cells.Item(10, 2).Value = 79
$ws.Cells.Item(11, 2).Value = 81
$ws.Cells.Item(12, 2).Value = 84
$ws.Cells.Item(13, 2).Value = 86
$ws.Cells.Item(14, 2).Value = 89
$ws.Cells.Item(15, 2).Value = 91
$ws.Cells.Item(16, 2).Value = 94
$ws.Cells.Item(17, 2).Value = 96
$ws.Cells.Item(18, 2).Value = 99
$ws.Cells.Item(19, 2).Value = 101
$ws.Cells.Item(20, 2).Value = 104
$ws.Cells.Item(21, 2).Value = 106
$ws.Cells.Item(22, 2).Value = 109
$ws.Cells.Item(23, 2).Value = 111
$ws.Cells.Item(24, 2).Value = 114
$ws.Cells.Item(25, 2).Value = 116
$ws.Cells.Item(26, 2).Value = 119
$ws.Cells.Item(27, 2).Value = 121
$ws.Cells.Item(28, 2).Value = 124
$ws.Cells.Item(29, 2).Value = 126
$ws.Cells.Item(30, 2).Value = 129
$ws.Cells.Item(31, 2).Value = 130
$ws.Cells.Item(32, 2).Value = 130
$ws.Cells.Item(33, 2).Value = 130
$ws.Cells.Item(34, 2).Value = 130

$ws = $wb.Worksheets.Item("2-Fall")
$ws.Cells.Item(7, 2).Value = 41
$ws.Cells.Item(8, 2).Value = 44
$ws.Cells.Item(9, 2).Value = 46
$ws.Cells.Item(10, 2).Value = 49
$ws.Cells.Item(11, 2).Value = 51
$ws.Cells.Item(12, 2).Value = 54
$ws.Cells.Item(13, 2).Value = 56
$ws.Cells.Item(14, 2).Value = 59
$ws.Cells.Item(15, 2).Value = 61
$ws.Cells.Item(16, 2).Value = 64
$ws.Cells.Item(17, 2).Value = 66
$ws.Cells.Item(18, 2).Value = 69
$ws.Cells.Item(19, 2).Value = 71
$ws.Cells.Item(20, 2).Value = 74
$ws.Cells.Item(21, 2).Value = 76
$ws.Cells.Item(22, 2).Value = 79
$ws.Cells.Item(23, 2).Value = 81
$ws.Cells.Item(24, 2).Value = 84
$ws.Cells.Item(25, 2).Value = 86
$ws.Cells.Item(26, 2).Value = 89
$ws.Cells.Item(27, 2).Value = 91
$ws.Cells.Item(28, 2).Value = 94
$ws.Cells.Item(29, 2).Value = 96
$ws.Cells.Item(30, 2).Value = 99
$ws.Cells.Item(31, 2).Value = 101
$ws.Cells.Item(32, 2).Value = 104
$ws.Cells.Item(33, 2).Value = 106
$ws.Cells.Item(34, 2).Value = 109

$ws = $wb.Worksheets.Item("2-Spring")
$ws.Cells.Item(10, 2).Value = 41
$ws.Cells.Item(11, 2).Value = 44
$ws.Cells.Item(12, 2).Value = 46
$ws.Cells.Item(13, 2).Value = 49
$ws.Cells.Item(14, 2).Value = 51
$ws.Cells.Item(15, 2).Value = 54
$ws.Cells.Item(16, 2).Value = 56
$ws.Cells.Item(17, 2).Value = 59
$ws.Cells.Item(18, 2).Value = 61
$ws.Cells.Item(19, 2).Value = 64
$ws.Cells.Item(20, 2).Value = 66
$ws.Cells.Item(21, 2).Value = 69
$ws.Cells.Item(22, 2).Value = 71
$ws.Cells.Item(23, 2).Value = 74
$ws.Cells.Item(24, 2).Value = 76
$ws.Cells.Item(25, 2).Value = 79
$ws.Cells.Item(26, 2).Value = 81
$ws.Cells.Item(27, 2).Value = 84
$ws.Cells.Item(28, 2).Value = 86
$ws.Cells.Item(29, 2).Value = 89
$ws.Cells.Item(30, 2).Value = 91
$ws.Cells.Item(31, 2).Value = 94
$ws.Cells.Item(32, 2).Value = 96
$ws.Cells.Item(33, 2).Value = 99
